# "future plot real updated"
# Swap the order of the first two slides: the slide that used to be
# second (the "学术/理论方法/双循环..." deck) should now come first,
# and the slide that used to be first (the "科研" deck) should move
# to the second position. The third slide is untouched.

$p = $ppt.ActivePresentation
$p.Slides.Item(2).MoveTo(1)
